$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New shared strings must be created in this order (B column first, then A
# column) so the resulting sharedStrings.xml ordering matches the target:
# CaoUser, Gemma Hardy, Admin User, Sahil Mittal.
$ws.Range("B1").Value = "CaoUser"
$ws.Range("B2").Value = "Gemma Hardy"
$ws.Range("A1").Value = "Admin User"
$ws.Range("A2").Value = "Sahil Mittal"

# Header row uses the bold style already applied to A1; apply the same to
# the new B1 header cell.
$ws.Range("B1").Font.Bold = $true

# Match (as closely as this engine's column-width model allows) the
# bestFit-derived widths from the real edit.
$ws.Columns("A").ColumnWidth = 10.6640625
$ws.Columns("B").ColumnWidth = 10.44140625

# Update the selection on the Users sheet and make it the active sheet/tab
# (the active tab moves from AddOpportunity back to Users).
$ws.Range("C11").Select() | Out-Null
$ws.Activate() | Out-Null
